# Auto-generated Excel COM-interop script
# Applies scheduled-runner profit recalculation updates across the Sagittarius_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1700.125
$ws.Range("J17").Value = 1700.125
$ws.Range("L17").Value = 5100.375
$ws.Range("N17").Value = -5436.375

$ws.Range("H41").Value = 188
$ws.Range("I41").Value = 188
$ws.Range("K41").Value = 188
$ws.Range("M41").Value = 252

$ws.Range("H70").Value = 1628.7
$ws.Range("I70").Value = 997
$ws.Range("K70").Value = 2991
$ws.Range("M70").Value = -2721

$ws.Range("H73").Value = 1628.7
$ws.Range("I73").Value = 997
$ws.Range("K73").Value = 2991
$ws.Range("M73").Value = -2055

$ws.Range("H99").Value = 255.6
$ws.Range("I99").Value = 265.75
$ws.Range("K99").Value = 797.25
$ws.Range("M99").Value = 700.75

$ws.Range("H101").Value = 33333842
$ws.Range("I101").Value = 50000264
$ws.Range("K101").Value = 150000792
$ws.Range("M101").Value = -149999170

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws.Range("H129").Value = 2650.9473
$ws.Range("I129").Value = 230
$ws.Range("J129").Value = 3104.875
$ws.Range("K129").Value = 690
$ws.Range("L129").Value = 9314.625
$ws.Range("M129").Value = 4310
$ws.Range("N129").Value = -19314.625

$ws.Range("H138").Value = 4702.423
$ws.Range("J138").Value = 5178.5
$ws.Range("L138").Value = 15535.5
$ws.Range("N138").Value = -25815.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2636.4666
$ws.Range("J61").Value = 3433.3333
$ws.Range("L61").Value = 3433.3333
$ws.Range("N61").Value = -3857.3333

$ws.Range("H74").Value = 582.7
$ws.Range("I74").Value = 582.7
$ws.Range("K74").Value = 582.7
$ws.Range("M74").Value = 291.3

$ws.Range("H77").Value = 582.7
$ws.Range("I77").Value = 582.7
$ws.Range("K77").Value = 2913.5
$ws.Range("M77").Value = 1454.5

$ws.Range("H102").Value = 1673
$ws.Range("I102").Value = 1673
$ws.Range("K102").Value = 1673
$ws.Range("M102").Value = -51

$ws.Range("H136").Value = 2636.4666
$ws.Range("J136").Value = 3433.3333
$ws.Range("L136").Value = 10299.9999
$ws.Range("N136").Value = -15399.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H132").Value = 2581.6
$ws.Range("I132").Value = 2816.2666
$ws.Range("J132").Value = 1877.6
$ws.Range("K132").Value = 8448.7998
$ws.Range("L132").Value = 5632.799999999999
$ws.Range("M132").Value = -5918.799800000001
$ws.Range("N132").Value = -10692.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 409.8421
$ws.Range("I44").Value = 504.4
$ws.Range("J44").Value = 376.07144
$ws.Range("K44").Value = 1513.2
$ws.Range("L44").Value = 1128.21432
$ws.Range("M44").Value = -1115.2
$ws.Range("N44").Value = -1924.21432

$ws.Range("H56").Value = 11024.378
$ws.Range("I56").Value = 11024.378
$ws.Range("K56").Value = 11024.378
$ws.Range("M56").Value = -10494.378

$ws.Range("H80").Value = 6400
$ws.Range("I80").Value = 4666.6665
$ws.Range("K80").Value = 13999.9995
$ws.Range("M80").Value = -13063.9995

$ws.Range("H83").Value = 6400
$ws.Range("I83").Value = 4666.6665
$ws.Range("K83").Value = 41999.9985
$ws.Range("M83").Value = -37319.9985

$ws.Range("H109").Value = 2350
$ws.Range("I109").Value = 562.5
$ws.Range("J109").Value = 4733.3335
$ws.Range("K109").Value = 1687.5
$ws.Range("L109").Value = 14200.0005
$ws.Range("M109").Value = -647.5
$ws.Range("N109").Value = -16280.0005

$ws.Range("H112").Value = 10300.333
$ws.Range("I112").Value = 2540.6
$ws.Range("J112").Value = 20000
$ws.Range("K112").Value = 7621.799999999999
$ws.Range("L112").Value = 60000
$ws.Range("M112").Value = -6513.799999999999
$ws.Range("N112").Value = -62216

$ws.Range("H115").Value = 3488.889
$ws.Range("I115").Value = 700
$ws.Range("J115").Value = 4285.7144
$ws.Range("K115").Value = 2100
$ws.Range("L115").Value = 12857.1432
$ws.Range("M115").Value = -925
$ws.Range("N115").Value = -15207.1432

$ws.Range("H141").Value = 7565.125
$ws.Range("I141").Value = 7565.125
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 22695.375
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -17515.375
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2365.2222
$ws.Range("I122").Value = 2057
$ws.Range("J122").Value = 2981.6667
$ws.Range("K122").Value = 6171
$ws.Range("L122").Value = 8945.000100000001
$ws.Range("M122").Value = -3721
$ws.Range("N122").Value = -13845.0001

$ws.Range("H140").Value = 143118.4
$ws.Range("J140").Value = 143118.4
$ws.Range("L140").Value = 143118.4
$ws.Range("N140").Value = -153478.4

$ws.Range("H141").Value = 59998
$ws.Range("J141").Value = 59998
$ws.Range("L141").Value = 59998
$ws.Range("N141").Value = -70358

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9999
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H122").Value = 6265.7407
$ws.Range("I122").Value = 4999.6665
$ws.Range("J122").Value = 7278.6
$ws.Range("K122").Value = 14998.9995
$ws.Range("L122").Value = 21835.8
$ws.Range("M122").Value = -12548.9995
$ws.Range("N122").Value = -26735.8

$ws.Range("H126").Value = 9999
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 6210.2
$ws.Range("I132").Value = 7354.4
$ws.Range("J132").Value = 3921.8
$ws.Range("K132").Value = 22063.2
$ws.Range("L132").Value = 11765.4
$ws.Range("M132").Value = -19533.2
$ws.Range("N132").Value = -16825.4

$ws.Range("H135").Value = 64949.5
$ws.Range("J135").Value = 64949.5
$ws.Range("L135").Value = 64949.5
$ws.Range("N135").Value = -75089.5

$ws.Range("H138").Value = 129990
$ws.Range("J138").Value = 129990
$ws.Range("L138").Value = 129990
$ws.Range("N138").Value = -140270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 19996.834
$ws.Range("I14").Value = 34992
$ws.Range("J14").Value = 12499.25
$ws.Range("K14").Value = 34992
$ws.Range("L14").Value = 12499.25
$ws.Range("M14").Value = -34824
$ws.Range("N14").Value = -12835.25

$ws.Range("H41").Value = 19978.666
$ws.Range("J41").Value = 19981.75
$ws.Range("L41").Value = 19981.75
$ws.Range("N41").Value = -20761.75

$ws.Range("H81").Value = 2000980.4
$ws.Range("I81").Value = 1000.3333
$ws.Range("J81").Value = 5000950.5
$ws.Range("K81").Value = 2000.6666
$ws.Range("L81").Value = 10001901
$ws.Range("M81").Value = -939.6666
$ws.Range("N81").Value = -10004023

$ws.Range("H84").Value = 2000980.4
$ws.Range("I84").Value = 1000.3333
$ws.Range("J84").Value = 5000950.5
$ws.Range("K84").Value = 10003.333
$ws.Range("L84").Value = 50009505
$ws.Range("M84").Value = -4699.333000000001
$ws.Range("N84").Value = -50020113

$ws.Range("H100").Value = 11112353
$ws.Range("I100").Value = 16667828
$ws.Range("J100").Value = 1403
$ws.Range("K100").Value = 33335656
$ws.Range("L100").Value = 2806
$ws.Range("M100").Value = -33335115
$ws.Range("N100").Value = -3888

$ws.Range("H107").Value = 594.5625
$ws.Range("I107").Value = 458.55554
$ws.Range("J107").Value = 769.4286
$ws.Range("K107").Value = 1375.66662
$ws.Range("L107").Value = 2308.2858
$ws.Range("M107").Value = 544.33338
$ws.Range("N107").Value = -6148.2858

$ws.Range("H122").Value = 1749
$ws.Range("I122").Value = 1749
$ws.Range("K122").Value = 5247
$ws.Range("M122").Value = -2797

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value = 34799
$ws.Range("J125").Value = 34799
$ws.Range("L125").Value = 34799
$ws.Range("N125").Value = -44639

$ws.Range("H126").Value = 5112.9375
$ws.Range("I126").Value = 4871.1763
$ws.Range("K126").Value = 14613.5289
$ws.Range("M126").Value = -12143.5289

$ws.Range("H132").Value = 2257.5
$ws.Range("I132").Value = 1850.4615
$ws.Range("J132").Value = 7549
$ws.Range("K132").Value = 5551.3845
$ws.Range("L132").Value = 22647
$ws.Range("M132").Value = -3021.3845
$ws.Range("N132").Value = -27707

$ws.Range("H141").Value = 206149
$ws.Range("J141").Value = 204579
$ws.Range("L141").Value = 204579
$ws.Range("N141").Value = -214939
